$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 8166.75
$ws.Range("I21").Value = 3000.3333
$ws.Range("K21").Value = 3000.3333
$ws.Range("M21").Value = -2532.3333
$ws.Range("H23").Value = 8166.75
$ws.Range("I23").Value = 3000.3333
$ws.Range("K23").Value = 3000.3333
$ws.Range("M23").Value = -2766.3333
$ws.Range("H31").Value = 857.8
$ws.Range("I31").Value = 929.6667
$ws.Range("J31").Value = 750
$ws.Range("K31").Value = 2789.0001
$ws.Range("L31").Value = 2250
$ws.Range("M31").Value = -2559.0001
$ws.Range("N31").Value = -2710
$ws.Range("H51").Value = 18717.084
$ws.Range("I51").Value = 30129.285
$ws.Range("J51").Value = 2740
$ws.Range("K51").Value = 30129.285
$ws.Range("L51").Value = 2740
$ws.Range("M51").Value = -29645.285
$ws.Range("N51").Value = -3708
$ws.Range("H106").Value = 2248.5
$ws.Range("I106").Value = 2093.25
$ws.Range("J106").Value = 3024.75
$ws.Range("K106").Value = 2093.25
$ws.Range("L106").Value = 3024.75
$ws.Range("M106").Value = -1462.25
$ws.Range("N106").Value = -4286.75
$ws.Range("H129").Value = 924.14545
$ws.Range("J129").Value = 965.0909
$ws.Range("L129").Value = 2895.2727
$ws.Range("N129").Value = -12895.2727

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4980252.5
$ws.Range("I32").Value = 4416.3857
$ws.Range("K32").Value = 4416.3857
$ws.Range("M32").Value = -4129.3857
$ws.Range("H45").Value = 18574908
$ws.Range("I45").Value = 33434032
$ws.Range("J45").Value = 1004.625
$ws.Range("K45").Value = 33434032
$ws.Range("L45").Value = 1004.625
$ws.Range("M45").Value = -33433655
$ws.Range("N45").Value = -1758.625
$ws.Range("H74").Value = 1901.8148
$ws.Range("I74").Value = 1379.1666
$ws.Range("J74").Value = 2319.9333
$ws.Range("K74").Value = 1379.1666
$ws.Range("L74").Value = 2319.9333
$ws.Range("M74").Value = -505.1666
$ws.Range("N74").Value = -4067.9333
$ws.Range("H77").Value = 1901.8148
$ws.Range("I77").Value = 1379.1666
$ws.Range("J77").Value = 2319.9333
$ws.Range("K77").Value = 6895.833000000001
$ws.Range("L77").Value = 11599.6665
$ws.Range("M77").Value = -2527.833000000001
$ws.Range("N77").Value = -20335.6665
$ws.Range("H102").Value = 1976.3914
$ws.Range("I102").Value = 1888.4286
$ws.Range("J102").Value = 2900
$ws.Range("K102").Value = 1888.4286
$ws.Range("L102").Value = 2900
$ws.Range("M102").Value = -266.4286
$ws.Range("N102").Value = -6144
$ws.Range("H109").Value = 46629.332
$ws.Range("J109").Value = 46629.332
$ws.Range("L109").Value = 46629.332
$ws.Range("N109").Value = -49403.332
$ws.Range("H133").Value = 29756.4
$ws.Range("J133").Value = 29756.4
$ws.Range("L133").Value = 29756.4
$ws.Range("N133").Value = -34816.4

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 993.37
$ws.Range("I105").Value = 992.8
$ws.Range("J105").Value = 995.08
$ws.Range("K105").Value = 992.8
$ws.Range("L105").Value = 995.08
$ws.Range("M105").Value = 754.2
$ws.Range("N105").Value = -4489.08
$ws.Range("H134").Value = 1737896.1
$ws.Range("I134").Value = 1544.3019
$ws.Range("K134").Value = 4632.905699999999
$ws.Range("M134").Value = -2097.905699999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 962.5231
$ws.Range("I31").Value = 699.82355
$ws.Range("J31").Value = 1250.6451
$ws.Range("K31").Value = 699.82355
$ws.Range("L31").Value = 1250.6451
$ws.Range("M31").Value = -404.82355
$ws.Range("N31").Value = -1840.6451
$ws.Range("H34").Value = 962.5231
$ws.Range("I34").Value = 699.82355
$ws.Range("J34").Value = 1250.6451
$ws.Range("K34").Value = 699.82355
$ws.Range("L34").Value = 1250.6451
$ws.Range("M34").Value = -497.82355
$ws.Range("N34").Value = -1654.6451
$ws.Range("H58").Value = 37038084
$ws.Range("J58").Value = 2725
$ws.Range("L58").Value = 2725
$ws.Range("N58").Value = -3131
$ws.Range("H122").Value = 20835356
$ws.Range("I122").Value = 25002134
$ws.Range("J122").Value = 1457
$ws.Range("K122").Value = 75006402
$ws.Range("L122").Value = 4371
$ws.Range("M122").Value = -75003952
$ws.Range("N122").Value = -9271
$ws.Range("H136").Value = 37038084
$ws.Range("J136").Value = 2725
$ws.Range("L136").Value = 8175
$ws.Range("N136").Value = -13275

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 608.84
$ws.Range("I107").Value = 171.19048
$ws.Range("J107").Value = 725.1772
$ws.Range("K107").Value = 513.5714400000001
$ws.Range("L107").Value = 2175.5316
$ws.Range("M107").Value = 1406.42856
$ws.Range("N107").Value = -6015.5316
$ws.Range("H118").Value = 23633.8
$ws.Range("I118").Value = 1686
$ws.Range("J118").Value = 56555.5
$ws.Range("K118").Value = 5058
$ws.Range("L118").Value = 169666.5
$ws.Range("M118").Value = -3815
$ws.Range("N118").Value = -172152.5
$ws.Range("H122").Value = 7249898
$ws.Range("I122").Value = 26315926
$ws.Range("J122").Value = 4807.12
$ws.Range("K122").Value = 236843334
$ws.Range("L122").Value = 43264.08
$ws.Range("M122").Value = -236840884
$ws.Range("N122").Value = -48164.08
$ws.Range("H131").Value = 882.35
$ws.Range("J131").Value = 908.54346
$ws.Range("L131").Value = 2725.63038
$ws.Range("N131").Value = -12805.63038

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1950
$ws.Range("I113").Value = 1675
$ws.Range("K113").Value = 1675
$ws.Range("M113").Value = 495
$ws.Range("H122").Value = 19294636
$ws.Range("I122").Value = 30012308
$ws.Range("J122").Value = 2829.1
$ws.Range("K122").Value = 90036924
$ws.Range("L122").Value = 8487.299999999999
$ws.Range("M122").Value = -90034474
$ws.Range("N122").Value = -13387.3

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1330.5
$ws.Range("I61").Value = 1330.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1330.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1128.5
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 1330.5
$ws.Range("I113").Value = 1330.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1330.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 839.5
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 7118.6113
$ws.Range("I122").Value = 8650.714
$ws.Range("K122").Value = 25952.142
$ws.Range("M122").Value = -23502.142
$ws.Range("H136").Value = 75572210
$ws.Range("I136").Value = 88437784
$ws.Range("J136").Value = 55559100
$ws.Range("K136").Value = 265313352
$ws.Range("L136").Value = 166677300
$ws.Range("M136").Value = -265310802
$ws.Range("N136").Value = -166682400

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 740.6
$ws.Range("I81").Value = 590.4
$ws.Range("J81").Value = 815.7
$ws.Range("K81").Value = 1180.8
$ws.Range("L81").Value = 1631.4
$ws.Range("M81").Value = -119.8
$ws.Range("N81").Value = -3753.4
$ws.Range("H84").Value = 740.6
$ws.Range("I84").Value = 590.4
$ws.Range("J84").Value = 815.7
$ws.Range("K84").Value = 5904
$ws.Range("L84").Value = 8157
$ws.Range("M84").Value = -600
$ws.Range("N84").Value = -18765
$ws.Range("H113").Value = 62503440
$ws.Range("I113").Value = 76926960
$ws.Range("K113").Value = 230780880
$ws.Range("M113").Value = -230778710
$ws.Range("H122").Value = 58099.223
$ws.Range("I122").Value = 100798.6
$ws.Range("K122").Value = 302395.8
$ws.Range("M122").Value = -299945.8
$ws.Range("H130").Value = 53133.332
$ws.Range("J130").Value = 53133.332
$ws.Range("L130").Value = 53133.332
$ws.Range("N130").Value = -63173.332
$ws.Range("H140").Value = 57420
$ws.Range("J140").Value = 57420
$ws.Range("L140").Value = 57420
$ws.Range("N140").Value = -67780
